$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the "TOPIC OF ALERT" text in C52 to include the topic path
# ------------------------------------------------------------------
$ws.Range("C52").Value = "TOPIC OF ALERT : Battery/IoT/project/UserID/1/statecontrol/AlertSMS "

# ------------------------------------------------------------------
# 2. Add the "done" markers in column D for rows 53-55, reusing the
#    same direct formatting already used by D39 (style applied to the
#    other "done" cells in the sheet).
# ------------------------------------------------------------------
$ws.Range("D39").Copy()
$ws.Range("D53:D55").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D53").Value = "done"
$ws.Range("D54").Value = "done"
$ws.Range("D55").Value = "done"

# ------------------------------------------------------------------
# 3. Add the new explanatory notes in column C for rows 53-55, reusing
#    the same direct formatting already used by C52.
#    NOTE: values are written in the same order they were appended to
#    the shared string table by the original author (casa / bassa /
#    alta) so that row 54 maps to "troppo alta" and row 55 maps to
#    "troppo bassa".
# ------------------------------------------------------------------
$ws.Range("C52").Copy()
$ws.Range("C53:C55").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("C53").Value = "se funziona, se non funziona, se la macchina è a casa"
$ws.Range("C55").Value = "se funziona, se non funziona, se è troppo bassa"
$ws.Range("C54").Value = "se troppo alta, non può caricare"

# ------------------------------------------------------------------
# 4. Update the sheet view: scroll down one row and move the
#    selection from C52 to D52.
# ------------------------------------------------------------------
$excel.Goto($ws.Range("A35"), $true) | Out-Null
$ws.Range("D52").Select() | Out-Null
